# "price feature" - insert a new "Price" column between the existing
# "ItemName" (A) and "First Name" (B) columns. The new column is left
# blank for now (dirty / unfilled price data), matching the commit
# message "prints out dirty data of price version".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift "First Name"/"Last Name" (and all the data beneath them) one
# column to the right, freeing up column B for the new header.
$ws.Columns("B").Insert()

# New header for the inserted column.
$ws.Range("B1").Value = "Price"

# Leave the active selection where the user was last working.
$ws.Range("C2").Select()
